$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.856.96'
$ws.Range('D3').Value = '2.223.32'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -2.27%  '
$ws.Range('D5').Value = '299.70'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').Value = '90.44'
$ws.Range('E6').Value = '  -3.88%  '
$ws.Range('D7').Value = '0.552'
$ws.Range('E7').Value = '  -3.33%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('E9').Value = '  -5.86%  '
$ws.Range('D10').Value = '33.05'
$ws.Range('E10').Value = '  -4.57%  '
$ws.Range('E11').Value = '  -3.31%  '
$ws.Range('E12').Value = '  -3.60%  '
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('D14').Value = '2.561.54'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '2.226.23'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '0.776'
$ws.Range('E17').Value = '  -6.95%  '
$ws.Range('D18').Value = '43.739.38'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').Value = '0.0₃0902'
$ws.Range('E19').Value = '  -5.42%  '
$ws.Range('D20').Value = '11.49'
$ws.Range('E20').Value = '  -3.91%  '
$ws.Range('E21').Value = '  -6.40%  '
$ws.Range('D22').Value = '64.48'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '236.41'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').Value = '2.81'
$ws.Range('E24').Value = '  -4.98%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('E26').Value = '  -5.14%  '
$ws.Range('D27').Value = '38.52'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('D29').Value = '9.35'
$ws.Range('E29').Value = '  -4.28%  '
$ws.Range('D30').Value = '153.16'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('E31').Value = '  -3.83%  '
$ws.Range('D32').Value = '5.42'
$ws.Range('E32').Value = '  -8.50%  '
$ws.Range('D33').Value = '0.0757'
$ws.Range('E33').Value = '  -5.06%  '
$ws.Range('D34').Value = '2.49'
$ws.Range('E34').Value = '  -5.55%  '
$ws.Range('E36').Value = '  -8.53%  '
$ws.Range('D37').Value = '0.103'
$ws.Range('E37').Value = '  -6.52%  '
$ws.Range('E38').Value = '  -5.24%  '
$ws.Range('D39').Value = '0.0299'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').Value = '3.19'
$ws.Range('E40').Value = '  -5.55%  '
$ws.Range('E41').Value = '  -3.73%  '
$ws.Range('D42').Value = '13.19'
$ws.Range('E42').Value = '  -8.47%  '
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D44').Value = '1.832.74'
$ws.Range('E44').Value = '  +3.72%  '
$ws.Range('E45').Value = '  +12.53%  '
$ws.Range('E46').Value = '  -5.73%  '
$ws.Range('D47').Value = '68.04'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = '94.37'
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('D49').Value = '72.95'
$ws.Range('E49').Value = '  -8.03%  '
$ws.Range('D50').Value = '13.91'
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('D51').Value = '7.74'
$ws.Range('E51').Value = '  -3.96%  '
